# Auto-generated COM-interop script
# 1) Turn the current '总计' sheet (sheet5, fund-holdings style) into '2022-Q1'
#    with the new per-fund holdings data.
# 2) Add a brand-new '总计' sheet at the end with the historical quarter summary
#    (now including the 2022-Q1 row).

$wb = $excel.ActiveWorkbook
$q1 = $wb.Worksheets.Item("总计")

# ---- Clear existing summary data from the sheet that becomes 2022-Q1 ----
$q1.Range("A1:D5").ClearContents()

# ---- Header row (style copied from the existing styled header cell) ----
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# B1:D1 already carry the header style; copy it across to the new E1:H1 cells
$q1.Range("B1:D1").Copy()
$q1.Range("E1:H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Per-fund holdings rows (2..19) ----
# Column A (index) keeps the bold/bordered style already present on A2:A5;
# new rows A6:A19 need that style copied across from A2.
# row 2: 161914
$q1.Range("A2").Value = 0
$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "161914"
$q1.Range("C2").NumberFormat = "@"
$q1.Range("C2").Value = "万家创业板2年定期开放混合A"
$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "15.74"
$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "95.20"
$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "9.34"
$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "1.4701"
$q1.Range("H2").Value = 1
# row 3: 360006
$q1.Range("A3").Value = 1
$q1.Range("B3").NumberFormat = "@"
$q1.Range("B3").Value = "360006"
$q1.Range("C3").NumberFormat = "@"
$q1.Range("C3").Value = "光大保德信新增长混合"
$q1.Range("D3").NumberFormat = "@"
$q1.Range("D3").Value = "21.71"
$q1.Range("E3").NumberFormat = "@"
$q1.Range("E3").Value = "88.07"
$q1.Range("F3").NumberFormat = "@"
$q1.Range("F3").Value = "4.18"
$q1.Range("G3").NumberFormat = "@"
$q1.Range("G3").Value = "0.9075"
$q1.Range("H3").Value = 5
# row 4: 003713
$q1.Range("A4").Value = 2
$q1.Range("B4").NumberFormat = "@"
$q1.Range("B4").Value = "003713"
$q1.Range("C4").NumberFormat = "@"
$q1.Range("C4").Value = "英大睿盛灵活配置混合A"
$q1.Range("D4").NumberFormat = "@"
$q1.Range("D4").Value = "5.99"
$q1.Range("E4").NumberFormat = "@"
$q1.Range("E4").Value = "87.42"
$q1.Range("F4").NumberFormat = "@"
$q1.Range("F4").Value = "6.18"
$q1.Range("G4").NumberFormat = "@"
$q1.Range("G4").Value = "0.3702"
$q1.Range("H4").Value = 4
# row 5: 161915
$q1.Range("A5").Value = 3
$q1.Range("B5").NumberFormat = "@"
$q1.Range("B5").Value = "161915"
$q1.Range("C5").NumberFormat = "@"
$q1.Range("C5").Value = "万家创业板2年定期开放混合C"
$q1.Range("D5").NumberFormat = "@"
$q1.Range("D5").Value = "2.36"
$q1.Range("E5").NumberFormat = "@"
$q1.Range("E5").Value = "95.20"
$q1.Range("F5").NumberFormat = "@"
$q1.Range("F5").Value = "9.34"
$q1.Range("G5").NumberFormat = "@"
$q1.Range("G5").Value = "0.2204"
$q1.Range("H5").Value = 1
# row 6: 010676
$q1.Range("A6").Value = 4
$q1.Range("B6").NumberFormat = "@"
$q1.Range("B6").Value = "010676"
$q1.Range("C6").NumberFormat = "@"
$q1.Range("C6").Value = "光大保德信新机遇混合"
$q1.Range("D6").NumberFormat = "@"
$q1.Range("D6").Value = "4.20"
$q1.Range("E6").NumberFormat = "@"
$q1.Range("E6").Value = "85.57"
$q1.Range("F6").NumberFormat = "@"
$q1.Range("F6").Value = "5.09"
$q1.Range("G6").NumberFormat = "@"
$q1.Range("G6").Value = "0.2138"
$q1.Range("H6").Value = 5
# row 7: 570005
$q1.Range("A7").Value = 5
$q1.Range("B7").NumberFormat = "@"
$q1.Range("B7").Value = "570005"
$q1.Range("C7").NumberFormat = "@"
$q1.Range("C7").Value = "诺德成长优势混合"
$q1.Range("D7").NumberFormat = "@"
$q1.Range("D7").Value = "4.30"
$q1.Range("E7").NumberFormat = "@"
$q1.Range("E7").Value = "62.89"
$q1.Range("F7").NumberFormat = "@"
$q1.Range("F7").Value = "4.20"
$q1.Range("G7").NumberFormat = "@"
$q1.Range("G7").Value = "0.1806"
$q1.Range("H7").Value = 6
# row 8: 003714
$q1.Range("A8").Value = 6
$q1.Range("B8").NumberFormat = "@"
$q1.Range("B8").Value = "003714"
$q1.Range("C8").NumberFormat = "@"
$q1.Range("C8").Value = "英大睿盛灵活配置混合C"
$q1.Range("D8").NumberFormat = "@"
$q1.Range("D8").Value = "2.40"
$q1.Range("E8").NumberFormat = "@"
$q1.Range("E8").Value = "87.42"
$q1.Range("F8").NumberFormat = "@"
$q1.Range("F8").Value = "6.18"
$q1.Range("G8").NumberFormat = "@"
$q1.Range("G8").Value = "0.1483"
$q1.Range("H8").Value = 4
# row 9: 012036
$q1.Range("A9").Value = 7
$q1.Range("B9").NumberFormat = "@"
$q1.Range("B9").Value = "012036"
$q1.Range("C9").NumberFormat = "@"
$q1.Range("C9").Value = "诺德兴远优选一年持有期混合型证券投资基金"
$q1.Range("D9").NumberFormat = "@"
$q1.Range("D9").Value = "2.75"
$q1.Range("E9").NumberFormat = "@"
$q1.Range("E9").Value = "52.19"
$q1.Range("F9").NumberFormat = "@"
$q1.Range("F9").Value = "4.15"
$q1.Range("G9").NumberFormat = "@"
$q1.Range("G9").Value = "0.1141"
$q1.Range("H9").Value = 4
# row 10: 003446
$q1.Range("A10").Value = 8
$q1.Range("B10").NumberFormat = "@"
$q1.Range("B10").Value = "003446"
$q1.Range("C10").NumberFormat = "@"
$q1.Range("C10").Value = "英大睿鑫灵活配置混合A"
$q1.Range("D10").NumberFormat = "@"
$q1.Range("D10").Value = "0.59"
$q1.Range("E10").NumberFormat = "@"
$q1.Range("E10").Value = "89.46"
$q1.Range("F10").NumberFormat = "@"
$q1.Range("F10").Value = "5.23"
$q1.Range("G10").NumberFormat = "@"
$q1.Range("G10").Value = "0.0309"
$q1.Range("H10").Value = 8
# row 11: 003447
$q1.Range("A11").Value = 9
$q1.Range("B11").NumberFormat = "@"
$q1.Range("B11").Value = "003447"
$q1.Range("C11").NumberFormat = "@"
$q1.Range("C11").Value = "英大睿鑫灵活配置混合C"
$q1.Range("D11").NumberFormat = "@"
$q1.Range("D11").Value = "0.51"
$q1.Range("E11").NumberFormat = "@"
$q1.Range("E11").Value = "89.46"
$q1.Range("F11").NumberFormat = "@"
$q1.Range("F11").Value = "5.23"
$q1.Range("G11").NumberFormat = "@"
$q1.Range("G11").Value = "0.0267"
$q1.Range("H11").Value = 8
# row 12: 011765
$q1.Range("A12").Value = 10
$q1.Range("B12").NumberFormat = "@"
$q1.Range("B12").Value = "011765"
$q1.Range("C12").NumberFormat = "@"
$q1.Range("C12").Value = "兴银高端制造混合A"
$q1.Range("D12").NumberFormat = "@"
$q1.Range("D12").Value = "1.01"
$q1.Range("E12").NumberFormat = "@"
$q1.Range("E12").Value = "93.23"
$q1.Range("F12").NumberFormat = "@"
$q1.Range("F12").Value = "2.60"
$q1.Range("G12").NumberFormat = "@"
$q1.Range("G12").Value = "0.0263"
$q1.Range("H12").Value = 9
# row 13: 003561
$q1.Range("A13").Value = 11
$q1.Range("B13").NumberFormat = "@"
$q1.Range("B13").Value = "003561"
$q1.Range("C13").NumberFormat = "@"
$q1.Range("C13").Value = "诺德成长精选灵活配置混合A"
$q1.Range("D13").NumberFormat = "@"
$q1.Range("D13").Value = "0.60"
$q1.Range("E13").NumberFormat = "@"
$q1.Range("E13").Value = "53.97"
$q1.Range("F13").NumberFormat = "@"
$q1.Range("F13").Value = "4.18"
$q1.Range("G13").NumberFormat = "@"
$q1.Range("G13").Value = "0.0251"
$q1.Range("H13").Value = 4
# row 14: 007152
$q1.Range("A14").Value = 12
$q1.Range("B14").NumberFormat = "@"
$q1.Range("B14").Value = "007152"
$q1.Range("C14").NumberFormat = "@"
$q1.Range("C14").Value = "诺德策略精选混合"
$q1.Range("D14").NumberFormat = "@"
$q1.Range("D14").Value = "0.42"
$q1.Range("E14").NumberFormat = "@"
$q1.Range("E14").Value = "93.32"
$q1.Range("F14").NumberFormat = "@"
$q1.Range("F14").Value = "5.01"
$q1.Range("G14").NumberFormat = "@"
$q1.Range("G14").Value = "0.0210"
$q1.Range("H14").Value = 7
# row 15: 007133
$q1.Range("A15").Value = 13
$q1.Range("B15").NumberFormat = "@"
$q1.Range("B15").Value = "007133"
$q1.Range("C15").NumberFormat = "@"
$q1.Range("C15").Value = "嘉实长青竞争优势股票A"
$q1.Range("D15").NumberFormat = "@"
$q1.Range("D15").Value = "0.25"
$q1.Range("E15").NumberFormat = "@"
$q1.Range("E15").Value = "90.07"
$q1.Range("F15").NumberFormat = "@"
$q1.Range("F15").Value = "4.97"
$q1.Range("G15").NumberFormat = "@"
$q1.Range("G15").Value = "0.0124"
$q1.Range("H15").Value = 5
# row 16: 011766
$q1.Range("A16").Value = 14
$q1.Range("B16").NumberFormat = "@"
$q1.Range("B16").Value = "011766"
$q1.Range("C16").NumberFormat = "@"
$q1.Range("C16").Value = "兴银高端制造混合C"
$q1.Range("D16").NumberFormat = "@"
$q1.Range("D16").Value = "0.39"
$q1.Range("E16").NumberFormat = "@"
$q1.Range("E16").Value = "93.23"
$q1.Range("F16").NumberFormat = "@"
$q1.Range("F16").Value = "2.60"
$q1.Range("G16").NumberFormat = "@"
$q1.Range("G16").Value = "0.0101"
$q1.Range("H16").Value = 9
# row 17: 001608
$q1.Range("A17").Value = 15
$q1.Range("B17").NumberFormat = "@"
$q1.Range("B17").Value = "001608"
$q1.Range("C17").NumberFormat = "@"
$q1.Range("C17").Value = "英大策略优选混合C"
$q1.Range("D17").NumberFormat = "@"
$q1.Range("D17").Value = "0.03"
$q1.Range("E17").NumberFormat = "@"
$q1.Range("E17").Value = "89.86"
$q1.Range("F17").NumberFormat = "@"
$q1.Range("F17").Value = "6.06"
$q1.Range("G17").NumberFormat = "@"
$q1.Range("G17").Value = "0.0018"
$q1.Range("H17").Value = 5
# row 18: 007134
$q1.Range("A18").Value = 16
$q1.Range("B18").NumberFormat = "@"
$q1.Range("B18").Value = "007134"
$q1.Range("C18").NumberFormat = "@"
$q1.Range("C18").Value = "嘉实长青竞争优势股票C"
$q1.Range("D18").NumberFormat = "@"
$q1.Range("D18").Value = "0.02"
$q1.Range("E18").NumberFormat = "@"
$q1.Range("E18").Value = "90.07"
$q1.Range("F18").NumberFormat = "@"
$q1.Range("F18").Value = "4.97"
$q1.Range("G18").NumberFormat = "@"
$q1.Range("G18").Value = "0.0010"
$q1.Range("H18").Value = 5
# row 19: 003562
$q1.Range("A19").Value = 17
$q1.Range("B19").NumberFormat = "@"
$q1.Range("B19").Value = "003562"
$q1.Range("C19").NumberFormat = "@"
$q1.Range("C19").Value = "诺德成长精选灵活配置混合C"
$q1.Range("D19").NumberFormat = "@"
$q1.Range("D19").Value = "0.00"
$q1.Range("E19").NumberFormat = "@"
$q1.Range("E19").Value = "53.97"
$q1.Range("F19").NumberFormat = "@"
$q1.Range("F19").Value = "4.18"
$q1.Range("G19").Value = 0
$q1.Range("H19").Value = 4

# Copy the A-column style (bold + border + centered) down onto the new rows
$q1.Range("A2").Copy()
$q1.Range("A6:A19").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Rename this sheet to 2022-Q1 ----
$q1.Name = "2022-Q1"

# ---- Create the new 总计 (totals) sheet at the end of the workbook ----
$newTotal = $wb.Worksheets.Add()
$newTotal.Name = "总计"

# Header row, matching the original totals-sheet layout
$newTotal.Range("B1").Value = "日期"
$newTotal.Range("C1").Value = "持有数量(只)"
$newTotal.Range("D1").Value = "持有市值(亿元)"

$newTotal.Range("A2").Value = 0
$newTotal.Range("B2").Value = "2022-Q1"
$newTotal.Range("C2").Value = 18
$newTotal.Range("D2").Value = 3.78
$newTotal.Range("A3").Value = 1
$newTotal.Range("B3").Value = "2021-Q4"
$newTotal.Range("C3").Value = 4
$newTotal.Range("D3").Value = 2.12
$newTotal.Range("A4").Value = 2
$newTotal.Range("B4").Value = "2021-Q3"
$newTotal.Range("C4").Value = 6
$newTotal.Range("D4").Value = 2.73
$newTotal.Range("A5").Value = 3
$newTotal.Range("B5").Value = "2021-Q2"
$newTotal.Range("C5").Value = 15
$newTotal.Range("D5").Value = 5.67
$newTotal.Range("A6").Value = 4
$newTotal.Range("B6").Value = "2021-Q1"
$newTotal.Range("C6").Value = 2
$newTotal.Range("D6").Value = 1.2

# Copy formatting from the equivalent cells of the 2022-Q1 sheet (which used to
# be this totals sheet, so its style 2 header/index cells are the right template)
$q1.Range("B1:D1").Copy()
$newTotal.Range("B1:D1").PasteSpecial(-4122)
$q1.Range("A2:A5").Copy()
$newTotal.Range("A2:A6").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---- Move the new 总计 sheet to the very end (after 2022-Q1) ----
$newTotal.Move($null, $wb.Worksheets.Item($wb.Worksheets.Count))

Write-Host "done"
